$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The worksheet is protected; unprotect it so the cells below can be edited,
# then re-apply protection once the updates are in place.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure banner.
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-14 for illustrative purposes only and are subject to change."

# That text contains an embedded line break, which makes the host recompute
# row 7's height; restore it to the default (no explicit custom height),
# matching the original layout.
$ws.Rows.Item(7).AutoFit()

# Refresh the weight / percent-change figures for the EFA, EEM and Total rows.
$ws.Range("D2").Value = 0.8460060160427807
$ws.Range("E2").Value = 0.001728395061728394

$ws.Range("D3").Value = 0.1539939839572192
$ws.Range("E3").Value = 0.002712967986977732

$ws.Range("E4").Value = 0.001880013368983802

# Restore sheet protection.
$ws.Protect()
